# Apply a global rename of model names that appear in the row-1 header
# labels of the form "<Judge> tag for <Target> Scenario No. <N>".
#
# Renames:
#   GPT-4o             -> GPT-5-mini
#   ChatGPT-4o         -> ChatGPT-5-mini
#   Gemini-2.5-Pro     -> Grok-4-Fast
#   Claude-3.7-Sonnet  -> Mistral-Small-24b-2501
#
# The rename is applied to both the "judge" and "target" model name
# occurrences in every header cell of row 1 (column B onward); column A
# ("Row_ID") and all data rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count

for ($c = 2; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Value

    if ($text) {
        $text = $text -replace 'Claude-3\.7-Sonnet', 'Mistral-Small-24b-2501'
        $text = $text -replace 'Gemini-2\.5-Pro', 'Grok-4-Fast'
        $text = $text -replace 'ChatGPT-4o', 'ChatGPT-5-mini'
        $text = $text -replace 'GPT-4o', 'GPT-5-mini'

        $cell.Value = $text
    }
}
